$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "DB Updated Date" value (B1) to the new date string.
# Temporarily mark the cell as Text so Excel doesn't auto-convert the
# "2021.02.02" literal into a date serial number, then drop the
# formatting again so the cell keeps the workbook's default style.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "2021.02.02"
$ws.Cells.ClearFormats()

# Update the "HUS, PPE Applied Year" value (B2) to the new year.
$ws.Range("B2").Value = 2021
